$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("G9").Value = 2.25
$ws.Range("I9").Value = 3.6
$ws.Range("M9").Value = 1.13
$ws.Range("N9").Value = 6
$ws.Range("Q9").Value = 2.88
$ws.Range("R9").Value = 1.4
$ws.Range("Y9").Value = 10
$ws.Range("AA9").Value = 23
$ws.Range("AS9").Value = 351
$ws.Range("BB9").Value = 451

# Row 10
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 3.1
$ws.Range("I10").Value = 4.2
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 4.75
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.63
$ws.Range("Q10").Value = 2.4
$ws.Range("R10").Value = 1.53
$ws.Range("S10").Value = 1.53
$ws.Range("T10").Value = 2.38
$ws.Range("U10").Value = 2.05
$ws.Range("V10").Value = 1.7
$ws.Range("W10").Value = 6
$ws.Range("X10").Value = 8.5
$ws.Range("Y10").Value = 9
$ws.Range("AB10").Value = 34
$ws.Range("AC10").Value = 7
$ws.Range("AE10").Value = 17
$ws.Range("AF10").Value = 67
$ws.Range("AG10").Value = 501
$ws.Range("AH10").Value = 9.5
$ws.Range("AI10").Value = 19
$ws.Range("AK10").Value = 41
$ws.Range("AM10").Value = 41
$ws.Range("AP10").Value = 26
$ws.Range("AR10").Value = 67
$ws.Range("AS10").Value = 251
$ws.Range("AT10").Value = 2.38
$ws.Range("AU10").Value = 9
$ws.Range("AV10").Value = 67
$ws.Range("AX10").Value = 23
$ws.Range("AY10").Value = 34
$ws.Range("AZ10").Value = 81
$ws.Range("BA10").Value = 126
$ws.Range("BB10").Value = 351

# Row 17
$ws.Range("G17").Value = 1.25
$ws.Range("H17").Value = 6
$ws.Range("J17").Value = 1.62
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 7.5
$ws.Range("M17").Value = 1.01
$ws.Range("N17").Value = 23
$ws.Range("O17").Value = 1.1
$ws.Range("P17").Value = 7
$ws.Range("Q17").Value = 1.33
$ws.Range("R17").Value = 3.4
$ws.Range("S17").Value = 1.17
$ws.Range("T17").Value = 4.33
$ws.Range("U17").Value = 1.67
$ws.Range("V17").Value = 2.1
$ws.Range("W17").Value = 11
$ws.Range("X17").Value = 8.5
$ws.Range("Y17").Value = 9.5
$ws.Range("Z17").Value = 9
$ws.Range("AC17").Value = 23
$ws.Range("AD17").Value = 12
$ws.Range("AF17").Value = 41
$ws.Range("AG17").Value = 151
$ws.Range("AL17").Value = 51
$ws.Range("AN17").Value = 3.6
$ws.Range("AO17").Value = 5.5
$ws.Range("AP17").Value = 13
$ws.Range("AR17").Value = 29
$ws.Range("AT17").Value = 4.33
$ws.Range("AU17").Value = 8.5
$ws.Range("AV17").Value = 41
